$d = $word.ActiveDocument

# --- Change 1: color the micro-SD bullet run red ---
$rng0 = $d.Content
$null = $rng0.Find.Execute("1x micro-SD- card via integrated HUB IC bridge for SPI")
$rng0.Font.Color = 255

# --- Change 2: color parts of the SD card UHS II bullet red ---
# First run: "1x SD card UHS II 250MB/s" -> find the SECOND occurrence in the doc
# (the first occurrence is a standalone bullet with no trailing description).
$rngFirstHit = $d.Content
$null = $rngFirstHit.Find.Execute("1x SD card UHS II 250MB/s")
$rngFirstHit.Collapse(0)

$rngSecondHit = $d.Range($rngFirstHit.End, $d.Content.End)
$null = $rngSecondHit.Find.Execute("1x SD card UHS II 250MB/s")
$rngSecondHit.Font.Color = 255

# Second run, split: the leading part up to (and including) "SD CARD " becomes red;
# the remainder ("controller. Recommended: ") stays default color.
$rngTail = $d.Range($rngSecondHit.End, $d.Content.End)
$null = $rngTail.Find.Execute(" – This will be interfaced to the HUB controller via a USB – SD CARD ")
$rngTail.Font.Color = 255
